$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task/goal rows appended below the existing tracker rows (1-12).
# Row 13: Task "newGoal"
$ws.Range("A13").Value = "Task"
$ws.Range("B13").Value = "newGoal"
$ws.Range("C13").NumberFormat = "dd/MM/yyyy"
$ws.Range("C13").Value = 44247.64841393519
$ws.Range("D13").Value = "1|1"
$ws.Range("E13").Value = "0|0"
$ws.Range("F13").Value = "'0%"
$ws.Range("F13").ClearFormats()
$ws.Range("G13").Value = "Just Started"

# Row 14: Goal "uguigiuhiuh"
$ws.Range("A14").Value = "Goal"
$ws.Range("B14").Value = "uguigiuhiuh"
$ws.Range("C14").NumberFormat = "dd/MM/yyyy"
$ws.Range("C14").Value = 44239.74582532408
$ws.Range("D14").Value = "0|0"
$ws.Range("E14").Value = "0|0"
$ws.Range("F14").Value = "'0%"
$ws.Range("F14").ClearFormats()
$ws.Range("G14").Value = "Just Started"

# Row 15: Task "guyguih9o"
$ws.Range("A15").Value = "Task"
$ws.Range("B15").Value = "guyguih9o"
$ws.Range("C15").NumberFormat = "dd/MM/yyyy"
$ws.Range("C15").Value = 44238.746113761576
$ws.Range("D15").Value = "2|2"
$ws.Range("E15").Value = "0|0"
$ws.Range("F15").Value = "'0%"
$ws.Range("F15").ClearFormats()
$ws.Range("G15").Value = "Just Started"
